$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note text in D9 (shared string for row 9)
$ws.Range("D9").Value = "more data formatting, finished framework for building studentgrades df"

# Update B9 value from 4 to 7
$ws.Range("B9").Value = 7

# Update selection to B10
$ws.Range("B10").Select()
